# Apply updated dSF (column F) values per repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -7
    3  = -5
    4  = 2
    5  = 3
    7  = -3
    8  = -3
    9  = -1
    10 = 1
    11 = -1
    12 = -5
    13 = -2
    14 = -8
    15 = -4
    16 = -2
    17 = -2
    18 = 10
    19 = 0
    20 = -1
    21 = -1
    22 = -1
    23 = 2
    26 = 1
    27 = 7
    28 = 3
    30 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
